$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 2386.5
$ws.Cells.Item(19, 9).Value = 2846
$ws.Cells.Item(19, 10).Value = 2080.1667
$ws.Cells.Item(19, 11).Value = 2846
$ws.Cells.Item(19, 12).Value = 2080.1667
$ws.Cells.Item(19, 13).Value = -2671
$ws.Cells.Item(19, 14).Value = -2430.1667
$ws.Cells.Item(38, 8).Value = 466
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 3607.5386
$ws.Cells.Item(40, 9).Value = 2320
$ws.Cells.Item(40, 11).Value = 2320
$ws.Cells.Item(40, 13).Value = -2145
$ws.Cells.Item(43, 8).Value = 3760.6155
$ws.Cells.Item(43, 10).Value = 4762.125
$ws.Cells.Item(43, 12).Value = 4762.125
$ws.Cells.Item(43, 14).Value = -4900.125
$ws.Cells.Item(58, 8).Value = 20
$ws.Cells.Item(58, 9).Value = 20
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 60
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).ClearContents()
$ws.Cells.Item(58, 14).Value = 90
$ws.Cells.Item(112, 8).Value = 2439.2856
$ws.Cells.Item(112, 10).Value = 2011.5385
$ws.Cells.Item(112, 12).Value = 6034.6155
$ws.Cells.Item(112, 14).Value = -8250.6155

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 8246017
$ws.Cells.Item(32, 9).Value = 1563389.9
$ws.Cells.Item(32, 11).Value = 1563389.9
$ws.Cells.Item(32, 13).Value = -1563102.9
$ws.Cells.Item(45, 8).Value = 10518.237
$ws.Cells.Item(45, 9).Value = 9439.179
$ws.Cells.Item(45, 10).Value = 13539.6
$ws.Cells.Item(45, 11).Value = 9439.179
$ws.Cells.Item(45, 12).Value = 13539.6
$ws.Cells.Item(45, 13).Value = -9062.179
$ws.Cells.Item(45, 14).Value = -14293.6
$ws.Cells.Item(80, 8).Value = 19329.166
$ws.Cells.Item(80, 10).Value = 19995.455
$ws.Cells.Item(80, 12).Value = 19995.455
$ws.Cells.Item(80, 14).Value = -21991.455
$ws.Cells.Item(83, 8).Value = 19329.166
$ws.Cells.Item(83, 10).Value = 19995.455
$ws.Cells.Item(83, 12).Value = 59986.36500000001
$ws.Cells.Item(83, 14).Value = -69970.36500000001
$ws.Cells.Item(110, 8).Value = 4579.516
$ws.Cells.Item(110, 9).Value = 4858.64
$ws.Cells.Item(110, 11).Value = 4858.64
$ws.Cells.Item(110, 13).Value = -2813.64

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(27, 8).Value = 30542
$ws.Cells.Item(27, 10).Value = 30542
$ws.Cells.Item(27, 12).Value = 30542
$ws.Cells.Item(27, 14).Value = -30926
$ws.Cells.Item(82, 8).Value = 20031.6
$ws.Cells.Item(82, 10).Value = 24996
$ws.Cells.Item(82, 12).Value = 24996
$ws.Cells.Item(82, 14).Value = -25762
$ws.Cells.Item(85, 8).Value = 20031.6
$ws.Cells.Item(85, 10).Value = 24996
$ws.Cells.Item(85, 12).Value = 24996
$ws.Cells.Item(85, 14).Value = -27648
$ws.Cells.Item(88, 8).Value = 82177.39999999999
$ws.Cells.Item(88, 10).Value = 82177.39999999999
$ws.Cells.Item(88, 12).Value = 82177.39999999999
$ws.Cells.Item(88, 14).Value = -82989.39999999999
$ws.Cells.Item(91, 8).Value = 82177.39999999999
$ws.Cells.Item(91, 10).Value = 82177.39999999999
$ws.Cells.Item(91, 12).Value = 82177.39999999999
$ws.Cells.Item(91, 14).Value = -84985.39999999999
$ws.Cells.Item(94, 8).Value = 2141
$ws.Cells.Item(94, 9).Value = 1934.7142
$ws.Cells.Item(94, 10).Value = 3224
$ws.Cells.Item(94, 11).Value = 1934.7142
$ws.Cells.Item(94, 12).Value = 3224
$ws.Cells.Item(94, 13).Value = -1483.7142
$ws.Cells.Item(94, 14).Value = -4126
$ws.Cells.Item(105, 8).Value = 12386.875
$ws.Cells.Item(105, 10).Value = 25643.25
$ws.Cells.Item(105, 12).Value = 25643.25
$ws.Cells.Item(105, 14).Value = -29137.25
$ws.Cells.Item(134, 8).Value = 3085.8572
$ws.Cells.Item(134, 9).Value = 2689.8708
$ws.Cells.Item(134, 11).Value = 8069.6124
$ws.Cells.Item(134, 13).Value = -5534.6124

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 2543.5833
$ws.Cells.Item(16, 9).Value = 1762.125
$ws.Cells.Item(16, 10).Value = 4106.5
$ws.Cells.Item(16, 11).Value = 1762.125
$ws.Cells.Item(16, 12).Value = 4106.5
$ws.Cells.Item(16, 13).Value = -1475.125
$ws.Cells.Item(16, 14).Value = -4680.5
$ws.Cells.Item(58, 8).Value = 992.6667
$ws.Cells.Item(58, 9).Value = 994.5
$ws.Cells.Item(58, 10).Value = 989
$ws.Cells.Item(58, 11).Value = 994.5
$ws.Cells.Item(58, 12).Value = 989
$ws.Cells.Item(58, 13).Value = -791.5
$ws.Cells.Item(58, 14).Value = -1395
$ws.Cells.Item(87, 8).Value = 25023.572
$ws.Cells.Item(87, 10).Value = 25023.572
$ws.Cells.Item(87, 12).Value = 25023.572
$ws.Cells.Item(87, 14).Value = -27395.572
$ws.Cells.Item(90, 8).Value = 25023.572
$ws.Cells.Item(90, 10).Value = 25023.572
$ws.Cells.Item(90, 12).Value = 75070.716
$ws.Cells.Item(90, 14).Value = -86926.716
$ws.Cells.Item(113, 8).Value = 2543.5833
$ws.Cells.Item(113, 9).Value = 1762.125
$ws.Cells.Item(113, 10).Value = 4106.5
$ws.Cells.Item(113, 11).Value = 1762.125
$ws.Cells.Item(113, 12).Value = 4106.5
$ws.Cells.Item(113, 13).Value = 407.875
$ws.Cells.Item(113, 14).Value = -8446.5
$ws.Cells.Item(132, 8).Value = 2662.6562
$ws.Cells.Item(132, 9).Value = 2444.5386
$ws.Cells.Item(132, 10).Value = 3607.8333
$ws.Cells.Item(132, 11).Value = 7333.6158
$ws.Cells.Item(132, 12).Value = 10823.4999
$ws.Cells.Item(132, 13).Value = -4803.6158
$ws.Cells.Item(132, 14).Value = -15883.4999
$ws.Cells.Item(134, 8).Value = 2260.9285
$ws.Cells.Item(134, 9).Value = 1814.32
$ws.Cells.Item(134, 11).Value = 5442.96
$ws.Cells.Item(134, 13).Value = -2907.96
$ws.Cells.Item(136, 8).Value = 992.6667
$ws.Cells.Item(136, 9).Value = 994.5
$ws.Cells.Item(136, 10).Value = 989
$ws.Cells.Item(136, 11).Value = 2983.5
$ws.Cells.Item(136, 12).Value = 2967
$ws.Cells.Item(136, 13).Value = -433.5
$ws.Cells.Item(136, 14).Value = -8067

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(68, 8).Value = 2209.0881
$ws.Cells.Item(68, 10).Value = 2144.5652
$ws.Cells.Item(68, 12).Value = 6433.6956
$ws.Cells.Item(68, 14).Value = -8055.6956
$ws.Cells.Item(71, 8).Value = 2209.0881
$ws.Cells.Item(71, 10).Value = 2144.5652
$ws.Cells.Item(71, 12).Value = 19301.0868
$ws.Cells.Item(71, 14).Value = -27413.0868
$ws.Cells.Item(122, 8).Value = 334.13333
$ws.Cells.Item(122, 10).Value = 428
$ws.Cells.Item(122, 12).Value = 3852
$ws.Cells.Item(122, 14).Value = -8752
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 13).ClearContents()

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(11, 8).Value = 23088492
$ws.Cells.Item(11, 9).Value = 2000266.6
$ws.Cells.Item(11, 10).Value = 33632610
$ws.Cells.Item(11, 11).Value = 2000266.6
$ws.Cells.Item(11, 12).Value = 33632610
$ws.Cells.Item(11, 13).Value = -2000127.6
$ws.Cells.Item(11, 14).Value = -33632888
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).ClearContents()
$ws.Cells.Item(12, 13).ClearContents()
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(14, 8).Value = 100000000
$ws.Cells.Item(14, 9).Value = 100000000
$ws.Cells.Item(14, 11).Value = 100000000
$ws.Cells.Item(14, 13).Value = -99999832
$ws.Cells.Item(97, 8).Value = 682.4286
$ws.Cells.Item(97, 9).Value = 705.6667
$ws.Cells.Item(97, 11).Value = 705.6667
$ws.Cells.Item(97, 13).Value = -209.6667
$ws.Cells.Item(122, 8).Value = 3537.8572
$ws.Cells.Item(122, 9).Value = 3303
$ws.Cells.Item(122, 11).Value = 9909
$ws.Cells.Item(122, 13).Value = -7459
$ws.Cells.Item(123, 8).Value = 72202.164
$ws.Cells.Item(123, 10).Value = 72202.164
$ws.Cells.Item(123, 12).Value = 72202.164
$ws.Cells.Item(123, 14).Value = -77102.164

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(122, 8).Value = 3480.8
$ws.Cells.Item(122, 9).Value = 3349.75
$ws.Cells.Item(122, 10).Value = 4005
$ws.Cells.Item(122, 11).Value = 10049.25
$ws.Cells.Item(122, 12).Value = 12015
$ws.Cells.Item(122, 13).Value = -7599.25
$ws.Cells.Item(122, 14).Value = -16915

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(19, 8).Value = 5000
$ws.Cells.Item(19, 10).Value = 5000
$ws.Cells.Item(19, 12).Value = 5000
$ws.Cells.Item(19, 14).Value = -5348
$ws.Cells.Item(81, 8).Value = 251250660
$ws.Cells.Item(81, 10).Value = 2000
$ws.Cells.Item(81, 12).Value = 4000
$ws.Cells.Item(81, 14).Value = -6122
$ws.Cells.Item(84, 8).Value = 251250660
$ws.Cells.Item(84, 10).Value = 2000
$ws.Cells.Item(84, 12).Value = 20000
$ws.Cells.Item(84, 14).Value = -30608
$ws.Cells.Item(132, 8).Value = 6916.231
$ws.Cells.Item(132, 9).Value = 7305.0835
$ws.Cells.Item(132, 10).Value = 2250
$ws.Cells.Item(132, 11).Value = 21915.2505
$ws.Cells.Item(132, 12).Value = 6750
$ws.Cells.Item(132, 13).Value = -19385.2505
$ws.Cells.Item(132, 14).Value = -11810
